$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 233
$ws.Range("B233").Value = 7559468
$ws.Range("F233").Value = 'Liverpool Montevideo'
$ws.Range("G233").Value = 'CA River Plate'
$ws.Range("H233").Value = 2
$ws.Range("I233").Value = 1
$ws.Range("J233").Value = 'H'
$ws.Range("K233").Value = 1.7
$ws.Range("L233").Value = 3
$ws.Range("M233").Value = 5.75
$ws.Range("N233").Value = 1.833
$ws.Range("O233").Value = 3.2
$ws.Range("P233").Value = 4.5
$ws.Range("Q233").Value = -0.5
$ws.Range("R233").Value = 1.925
$ws.Range("S233").Value = 1.925
$ws.Range("T233").Value = 2.25
$ws.Range("U233").Value = 2.025
$ws.Range("V233").Value = 1.825
$ws.Range("W233").Value = 0.833
$ws.Range("X233").Value = -1
$ws.Range("Y233").Value = -1
$ws.Range("Z233").Value = 0.925
$ws.Range("AA233").Value = -1
$ws.Range("AB233").Value = 1.025
$ws.Range("AC233").Value = -1

# Row 234
$ws.Range("B234").Value = 7559469
$ws.Range("F234").Value = 'Montevideo Wanderers'
$ws.Range("G234").Value = 'Penarol'
$ws.Range("H234").Value = 0
$ws.Range("I234").Value = 0
$ws.Range("J234").Value = 'D'
$ws.Range("K234").Value = 4.75
$ws.Range("L234").Value = 3.4
$ws.Range("M234").Value = 1.7
$ws.Range("N234").Value = 2.7
$ws.Range("O234").Value = 3.2
$ws.Range("P234").Value = 2.45
$ws.Range("Q234").Value = 0
$ws.Range("R234").Value = 2.05
$ws.Range("S234").Value = 1.8
$ws.Range("T234").Value = 2.5
$ws.Range("U234").Value = 1.975
$ws.Range("V234").Value = 1.875
$ws.Range("W234").Value = -1
$ws.Range("X234").Value = 2.2
$ws.Range("Y234").Value = -1
$ws.Range("Z234").Value = 0
$ws.Range("AA234").Value = -0
$ws.Range("AB234").Value = -1
$ws.Range("AC234").Value = 0.875

# Row 236
$ws.Range("B236").Value = 7013409
$ws.Range("F236").Value = 'Nacional De Football'
$ws.Range("G236").Value = 'Torque'
$ws.Range("H236").Value = 1
$ws.Range("I236").Value = 1
$ws.Range("J236").Value = 'D'
$ws.Range("K236").Value = 1.666
$ws.Range("L236").Value = 3.9
$ws.Range("M236").Value = 4.5
$ws.Range("N236").Value = 1.615
$ws.Range("O236").Value = 4
$ws.Range("P236").Value = 4.75
$ws.Range("Q236").Value = -0.75
$ws.Range("R236").Value = 1.8
$ws.Range("S236").Value = 2.05
$ws.Range("T236").Value = 2.75
$ws.Range("U236").Value = 1.95
$ws.Range("V236").Value = 1.9
$ws.Range("W236").Value = -1
$ws.Range("X236").Value = 3
$ws.Range("Y236").Value = -1
$ws.Range("Z236").Value = -1
$ws.Range("AA236").Value = 1.05
$ws.Range("AB236").Value = -1
$ws.Range("AC236").Value = 0.8999999999999999

# Row 237
$ws.Range("B237").Value = 7013702
$ws.Range("F237").Value = 'Defensor Sporting'
$ws.Range("G237").Value = 'Danubio'
$ws.Range("H237").Value = 0
$ws.Range("I237").Value = 2
$ws.Range("J237").Value = 'A'
$ws.Range("K237").Value = 1.8
$ws.Range("L237").Value = 3.6
$ws.Range("M237").Value = 4.2
$ws.Range("N237").Value = 1.8
$ws.Range("O237").Value = 3.6
$ws.Range("P237").Value = 4.2
$ws.Range("Q237").Value = -0.75
$ws.Range("R237").Value = 2.05
$ws.Range("S237").Value = 1.8
$ws.Range("T237").Value = 2.25
$ws.Range("U237").Value = 1.85
$ws.Range("V237").Value = 2
$ws.Range("W237").Value = -1
$ws.Range("X237").Value = -1
$ws.Range("Y237").Value = 3.2
$ws.Range("Z237").Value = -1
$ws.Range("AA237").Value = 0.8
$ws.Range("AB237").Value = -0.5
$ws.Range("AC237").Value = 0.5

# Row 238
$ws.Range("B238").Value = 7013885
$ws.Range("F238").Value = 'La Luz'
$ws.Range("G238").Value = 'Atletico Fenix Montevideo'
$ws.Range("H238").Value = 0
$ws.Range("I238").Value = 2
$ws.Range("J238").Value = 'A'
$ws.Range("K238").Value = 3
$ws.Range("L238").Value = 3
$ws.Range("M238").Value = 2.4
$ws.Range("N238").Value = 2.9
$ws.Range("O238").Value = 2.75
$ws.Range("P238").Value = 2.6
$ws.Range("Q238").Value = 0
$ws.Range("R238").Value = 2.025
$ws.Range("S238").Value = 1.825
$ws.Range("T238").Value = 2
$ws.Range("U238").Value = 2.025
$ws.Range("V238").Value = 1.825
$ws.Range("W238").Value = -1
$ws.Range("X238").Value = -1
$ws.Range("Y238").Value = 1.6
$ws.Range("Z238").Value = -1
$ws.Range("AA238").Value = 0.825
$ws.Range("AB238").Value = 0
$ws.Range("AC238").Value = -0

# Row 239
$ws.Range("B239").Value = 7013886
$ws.Range("F239").Value = 'Racing Club de Montevideo'
$ws.Range("G239").Value = 'Cerro'
$ws.Range("H239").Value = 0
$ws.Range("I239").Value = 1
$ws.Range("J239").Value = 'A'
$ws.Range("K239").Value = 2.25
$ws.Range("L239").Value = 3.1
$ws.Range("M239").Value = 3.25
$ws.Range("N239").Value = 2.25
$ws.Range("O239").Value = 2.875
$ws.Range("P239").Value = 3.5
$ws.Range("Q239").Value = -0.25
$ws.Range("R239").Value = 1.95
$ws.Range("S239").Value = 1.9
$ws.Range("T239").Value = 2
$ws.Range("U239").Value = 1.925
$ws.Range("V239").Value = 1.925
$ws.Range("W239").Value = -1
$ws.Range("X239").Value = -1
$ws.Range("Y239").Value = 2.5
$ws.Range("Z239").Value = -1
$ws.Range("AA239").Value = 0.8999999999999999
$ws.Range("AB239").Value = -1
$ws.Range("AC239").Value = 0.925
